# The workbook's reading/plotting code was automated to pull ALL columns
# (including the security/product names) straight out of Excel, so this
# edit refreshes both the label column (A) and the recomputed optimizer
# output columns (B/C/D), and appends one new holding row (row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: security / product names --------------------------------
# New row 2 label ("thanos") plus every other label shifts to the next
# instrument in the updated universe; a new ninth holding is appended.
$ws.Range("A2").Value = "thanos"
$ws.Range("A3").Value = "PARSTEI LX Equity"
$ws.Range("A4").Value = "LEF1TREU Index"
$ws.Range("A5").Value = "SX5R Index"
$ws.Range("A6").Value = "SXUSR Index"
$ws.Range("A7").Value = "BEGCGA Index"
$ws.Range("A8").Value = "LEC4TREU Index"
$ws.Range("A9").Value = "LEATTREU Index"

# Give the new label cell (A9) the same header/label style used by the
# other rows in column A (bold, centered, bordered - style index 1).
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Columns B/C/D: weights + recomputed optimizer outputs -------------
$ws.Range("B2").Value = 0.3
$ws.Range("C2").Value = 0.03595218056188074
$ws.Range("D2").Value = 0.0359521717459174

$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = 0.1909876858665356
$ws.Range("D3").Value = 0.1909876982079544

$ws.Range("B4").Value = 0.15
$ws.Range("C4").Value = 0.1838869452973091
$ws.Range("D4").Value = 0.1838869262343344

$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = 0.07553560992936761
$ws.Range("D5").Value = 0.0755356137722062

$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = 0.06451041465115818
$ws.Range("D6").Value = 0.0645103867352906

$ws.Range("B7").Value = 0.05
$ws.Range("C7").Value = 0.1114101465131827
$ws.Range("D7").Value = 0.111410190157718

$ws.Range("B8").Value = 0.05
$ws.Range("C8").Value = 0.1699542531985195
$ws.Range("D8").Value = 0.1699542547564113

# --- New row 9: additional holding --------------------------------------
$ws.Range("B9").Value = 0.05
$ws.Range("C9").Value = 0.1677627639820466
$ws.Range("D9").Value = 0.1677627583901678
